# Remove the trailing "Ver no Jupiter..." and copyright paragraphs (and
# the blank paragraph that separated them from the Requisitos section),
# leaving the "LOB1012: Estatística (Requisito)" paragraph followed
# directly by the existing blank paragraph / page-break paragraph.

$d = $word.ActiveDocument

$target1 = "Ver no Jupiter Salvar em pdf Salvar em docx"
$target2 = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

# Locate the "Ver no Jupiter..." paragraph by its text and walk from there.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq $target1) {
        $jupiterPara = $p
        break
    }
}

$copyrightPara = $jupiterPara.Next()
$blankPara = $jupiterPara.Previous()

# Sanity checks on neighbouring text before mutating anything.
if ($copyrightPara.Range.Text.TrimEnd("`r", "`a") -ne $target2) {
    throw "Unexpected paragraph after 'Ver no Jupiter...' paragraph"
}
if ($blankPara.Range.Text.TrimEnd("`r", "`a") -ne "") {
    throw "Unexpected paragraph before 'Ver no Jupiter...' paragraph"
}

# Build one contiguous range spanning the blank paragraph through the
# copyright paragraph (inclusive) and delete it in one shot.
$deleteRange = $d.Range($blankPara.Range.Start, $copyrightPara.Range.End)
$deleteRange.Delete()

Write-Output "Paragraphs after edit:"
Write-Output $d.Paragraphs.Count
